# Applies the "Update gh-pages to output generated at 456a3b4" change:
# bumps the "想去人数" (column F) counts on the 展览, 演出 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 33
$ws.Range("F4").Value = 202
$ws.Range("F5").Value = 1112
$ws.Range("F6").Value = 8307
$ws.Range("F7").Value = 8307
$ws.Range("F9").Value = 214
$ws.Range("F10").Value = 6938
$ws.Range("F11").Value = 176
$ws.Range("F12").Value = 5076
$ws.Range("F13").Value = 5551
$ws.Range("F14").Value = 1080
$ws.Range("F15").Value = 343
$ws.Range("F16").Value = 350
$ws.Range("F17").Value = 26
$ws.Range("F19").Value = 257
$ws.Range("F20").Value = 137
$ws.Range("F23").Value = 152
$ws.Range("F25").Value = 9322
$ws.Range("F26").Value = 75
$ws.Range("F27").Value = 1712
$ws.Range("F28").Value = 1020
$ws.Range("F31").Value = 1901
$ws.Range("F37").Value = 1912
$ws.Range("F39").Value = 1220
$ws.Range("F41").Value = 4867
$ws.Range("F42").Value = 380
$ws.Range("F46").Value = 154
$ws.Range("F47").Value = 1083
$ws.Range("F48").Value = 1048
$ws.Range("F49").Value = 929
$ws.Range("F50").Value = 1277

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 36

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 33
$ws.Range("F4").Value = 202
$ws.Range("F6").Value = 1112
$ws.Range("F7").Value = 8307
$ws.Range("F9").Value = 214
$ws.Range("F10").Value = 6938
$ws.Range("F11").Value = 176
$ws.Range("F14").Value = 5077
$ws.Range("F15").Value = 5551
$ws.Range("F16").Value = 1080
$ws.Range("F17").Value = 343
$ws.Range("F18").Value = 350
$ws.Range("F19").Value = 26
$ws.Range("F21").Value = 257
$ws.Range("F22").Value = 137
$ws.Range("F23").Value = 152
$ws.Range("F25").Value = 9322
$ws.Range("F26").Value = 75
$ws.Range("F27").Value = 1712
$ws.Range("F28").Value = 1020
$ws.Range("F31").Value = 1901
$ws.Range("F37").Value = 1912
$ws.Range("F39").Value = 1220
$ws.Range("F41").Value = 4867
$ws.Range("F42").Value = 380
$ws.Range("F46").Value = 154
$ws.Range("F47").Value = 1083
$ws.Range("F48").Value = 1048
$ws.Range("F49").Value = 929
$ws.Range("F50").Value = 1277

